$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.021.23"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "'1.633.00"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'214.45"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").Value = "'0.502"
$ws.Range("E6").Value = "  -1.39%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -2.64%  "
$ws.Range("D9").Value = "'0.0619"
$ws.Range("E9").Value = "  -3.56%  "
$ws.Range("D10").Value = "'18.12"
$ws.Range("E10").Value = "  -7.97%  "
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").Value = "'1.859.17"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "'1.626.99"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("E14").Value = "  -3.27%  "
$ws.Range("E15").Value = "  -3.96%  "
$ws.Range("D16").Value = "'26.005.06"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "'0.0₃0739"
$ws.Range("E17").Value = "  -3.77%  "
$ws.Range("D18").Value = "'61.30"
$ws.Range("E18").Value = "  -3.60%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "'190.22"
$ws.Range("E20").Value = "  -3.22%  "
$ws.Range("D21").Value = "'4.23"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").Value = "'6.07"
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("D24").Value = "'0.131"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'143.86"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'6.75"
$ws.Range("E28").Value = "  -2.46%  "
$ws.Range("D29").Value = "'15.13"
$ws.Range("E29").Value = "  -3.26%  "
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("E31").Value = "  -3.79%  "
$ws.Range("D32").Value = "'3.12"
$ws.Range("E32").Value = "  -4.68%  "
$ws.Range("E33").Value = "  -5.73%  "
$ws.Range("D34").Value = "'2.42"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").Value = "'1.48"
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("D36").Value = "'1.125.45"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").Value = "'0.858"
$ws.Range("D38").Value = "'2.44"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("E39").Value = "  -4.99%  "
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("D41").Value = "'98.05"
$ws.Range("E41").Value = "  -1.53%  "
$ws.Range("D42").Value = "'0.775"
$ws.Range("E42").Value = "  -3.05%  "
$ws.Range("D43").Value = "'1.771.30"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("E44").Value = "  -5.35%  "
$ws.Range("D45").Value = "'0.0₆0111"
$ws.Range("E45").Value = "  -3.76%  "
$ws.Range("D46").Value = "'54.73"
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("D47").Value = "'0.0527"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "'1.48"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "'7.49"
$ws.Range("E51").Value = "  -3.37%  "
